$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: paragraph 1 ("los usuarios que estan suspendidos...") -> red text
# ---------------------------------------------------------------------------
$p1 = $d.Paragraphs(1)
$p1.Range.Font.Color = 255

# ---------------------------------------------------------------------------
# Change 2: paragraph 5 ("La barra de herramientas...") -> merge the two runs
# that were split by a _GoBack bookmark into a single run / remove bookmark
# ---------------------------------------------------------------------------
$p5 = $d.Paragraphs(5)
$null = $p5.Range.Find.Execute(
    "La barra de herramientas no muestra la foto del usuario",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "La barra de herramientas no muestra la foto del usuario", 2)

# ---------------------------------------------------------------------------
# Change 3: paragraph 9 ("capturar error cuando verificar el email esta
# duplicado") -> red text
# ---------------------------------------------------------------------------
$p9 = $d.Paragraphs(9)
$p9.Range.Font.Color = 255

# ---------------------------------------------------------------------------
# Change 4: paragraph 11 ("ordenar las tablas...") -> split "ordenar las t" /
# "ablas con un pipe y ponerle un " with a _GoBack bookmark in between
# ---------------------------------------------------------------------------
$p11 = $d.Paragraphs(11)
$p11Start = $p11.Range.Start
$bmRange = $d.Range($p11Start + 13, $p11Start + 13)
$null = $d.Bookmarks.Add("_GoBack", $bmRange)

# ---------------------------------------------------------------------------
# Change 5: paragraph 23 ("arreglar el ver encuestas...") -> fix the
# "qu emostrar" typo to "que mostrar" and restructure the runs / drop the
# now-stale spellcheck proofErr markers
# ---------------------------------------------------------------------------
$p23 = $d.Paragraphs(23)
$p23Body = $d.Range($p23.Range.Start, $p23.Range.End - 1)
$xml23 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">arreglar el ver encuestas, muestra todos los comentarios buenos y malos, tiene </w:t></w:r><w:r><w:t>qu</w:t></w:r><w:r><w:t>e</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>mostrar los mejores y peores puntuados</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$null = $p23Body.InsertXML($xml23)

# ---------------------------------------------------------------------------
# Change 6: paragraph 25 ("que no traiga todos los pedidos del cliente
# cuando los lista") -> "no muestra la imagen del pedido cuando los lista"
# ---------------------------------------------------------------------------
$p25 = $d.Paragraphs(25)
$null = $p25.Range.Find.Execute(
    "que no traiga todos los pedidos del cliente cuando los lista",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "no muestra la imagen del pedido cuando los lista", 2)

# ---------------------------------------------------------------------------
# Change 7: add a new bullet paragraph right before the trailing empty
# paragraph, after "... para cambiar la contraseña"
# ---------------------------------------------------------------------------
$pLast = $d.Paragraphs($d.Paragraphs.Count)
$insertPoint = $d.Range($pLast.Range.Start, $pLast.Range.Start)
$null = $insertPoint.InsertParagraphBefore()

$pNew = $d.Paragraphs($d.Paragraphs.Count - 1)
$pNewRange = $d.Range($pNew.Range.Start, $pNew.Range.End)
$xmlNew = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">cuando te </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>logueas</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> que espere para que se ponga en blanco todo</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$null = $pNewRange.InsertXML($xmlNew)

Write-Output "edit complete"
